$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.609.47"
$ws.Range("E2").Value = "  +2.22%  "

$ws.Range("D3").Value = "3.830.96"
$ws.Range("E3").Value = "  +0.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "631.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").Value = "3.834.69"
$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.455"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "

$ws.Range("D15").Value = "4.465.44"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").Value = "3.796.93"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "69.497.64"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.711"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "3.974.10"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").Value = "3.771.86"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("E38").Value = "  +3.76%  "

$ws.Range("E39").Value = "  +7.23%  "

$ws.Range("E40").Value = "  +5.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.981"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.38%  "
